$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.418.62"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").Value = "1.952.76"
$ws.Range("E3").Value = "  -1.66%  "

# Row 5
$ws.Range("D5").Value = "'244.19"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6
$ws.Range("E6").Value = "  -1.98%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'57.62"
$ws.Range("E8").Value = "  -1.93%  "

# Row 9
$ws.Range("D9").Value = "'0.366"
$ws.Range("E9").Value = "  -3.07%  "

# Row 10
$ws.Range("D10").Value = "'0.0853"
$ws.Range("E10").Value = "  +4.56%  "

# Row 11
$ws.Range("E11").Value = "  +0.43%  "

# Row 12
$ws.Range("D12").Value = "2.239.02"
$ws.Range("E12").Value = "  -1.54%  "

# Row 13
$ws.Range("E13").Value = "  -5.11%  "

# Row 14
$ws.Range("D14").Value = "'21.51"
$ws.Range("E14").Value = "  -10.97%  "

# Row 15
$ws.Range("D15").Value = "'13.56"
$ws.Range("E15").Value = "  -3.90%  "

# Row 16
$ws.Range("D16").Value = "'5.21"
$ws.Range("E16").Value = "  -4.67%  "

# Row 17
$ws.Range("D17").Value = "1.952.04"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18
$ws.Range("D18").Value = "36.352.59"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0886"
$ws.Range("E19").Value = "  +2.36%  "

# Row 20
$ws.Range("D20").Value = "'69.79"
$ws.Range("E20").Value = "  -1.89%  "

# Row 21
$ws.Range("D21").Value = "'230.15"
$ws.Range("E21").Value = "  -2.12%  "

# Row 22
$ws.Range("D22").Value = "'5.07"
$ws.Range("E22").Value = "  -4.98%  "

# Row 23
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("E24").Value = "  -7.32%  "

# Row 25
$ws.Range("E25").Value = "  -0.56%  "

# Row 26
$ws.Range("D26").Value = "'9.29"
$ws.Range("E26").Value = "  -9.12%  "

# Row 27
$ws.Range("D27").Value = "'161.63"
$ws.Range("E27").Value = "  -0.39%  "

# Row 28
$ws.Range("D28").Value = "'0.135"
$ws.Range("E28").Value = "  +7.87%  "

# Row 29
$ws.Range("D29").Value = "'19.61"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("E30").Value = "  -1.84%  "

# Row 31
$ws.Range("E31").Value = "  -2.86%  "

# Row 32
$ws.Range("E32").Value = "  -5.43%  "

# Row 33
$ws.Range("D33").Value = "'0.0649"
$ws.Range("E33").Value = "  +2.60%  "

# Row 34
$ws.Range("E34").Value = "  -3.96%  "

# Row 35
$ws.Range("E35").Value = "  -1.54%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("E37").Value = "  +1.15%  "

# Row 38
$ws.Range("E38").Value = "  -5.66%  "

# Row 39
$ws.Range("D39").Value = "'3.05"
$ws.Range("E39").Value = "  -1.50%  "

# Row 40
$ws.Range("D40").Value = "'0.0981"
$ws.Range("E40").Value = "  +1.14%  "

# Row 41
$ws.Range("E41").Value = "  +0.28%  "

# Row 42
$ws.Range("D42").Value = "'0.0213"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43
$ws.Range("E43").Value = "  -6.78%  "

# Row 44
$ws.Range("D44").Value = "'15.74"
$ws.Range("E44").Value = "  -3.72%  "

# Row 45
$ws.Range("D45").Value = "1.360.35"
$ws.Range("E45").Value = "  -0.87%  "

# Row 46
$ws.Range("E46").Value = "  -6.23%  "

# Row 47
$ws.Range("D47").Value = "'87.85"
$ws.Range("E47").Value = "  -5.64%  "

# Row 48
$ws.Range("E48").Value = "  -6.52%  "

# Row 49
$ws.Range("E49").Value = "  -0.51%  "

# Row 50
$ws.Range("D50").Value = "'45.09"
$ws.Range("E50").Value = "  -0.58%  "

# Row 51
$ws.Range("D51").Value = "2.129.68"
$ws.Range("E51").Value = "  -1.90%  "
